$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new group of 3 rows (one reporting date's worth of data) right
# before the current row 338. Excel shifts rows 338:442 down to 341:445
# (carrying formats), and the previously-last group (old rows 440:442,
# date 2021-12-01 / 44540) ends up re-numbered as new rows 443:445,
# unchanged. The freshly inserted rows 338:340 get populated with the new
# week's data (date 2022-01-01 / 44551).
$ws.Range("A338:A340").EntireRow.Insert()

# Row 338 - Pintón
$ws.Range("A338").Value = 8
$ws.Range("B338").Value = "Terminal La Palmera de La Serena"
$ws.Range("C338").Value = "Coquimbo"
$ws.Range("D338").Value = 44551
$ws.Range("E338").Value = 4
$ws.Range("F338").Value = "Fruta"
$ws.Range("G338").Value = 100108
$ws.Range("H338").Value = "Tropicales y subtropicales"
$ws.Range("I338").Value = 100108006
$ws.Range("J338").Value = "Plátano"
$ws.Range("K338").Value = "Sin especificar"
$ws.Range("L338").Value = "Pintón"
$ws.Range("M338").Value = 120
$ws.Range("N338").Value = 13500
$ws.Range("O338").Value = 13500
$ws.Range("P338").Value = 13500
$ws.Range("Q338").Value = "`$/caja 20 kilos"
$ws.Range("R338").Value = "Ecuador"
$ws.Range("S338").Value = 675
$ws.Range("T338").Value = 20

# Row 339 - Primera Maduro
$ws.Range("A339").Value = 8
$ws.Range("B339").Value = "Terminal La Palmera de La Serena"
$ws.Range("C339").Value = "Coquimbo"
$ws.Range("D339").Value = 44551
$ws.Range("E339").Value = 4
$ws.Range("F339").Value = "Fruta"
$ws.Range("G339").Value = 100108
$ws.Range("H339").Value = "Tropicales y subtropicales"
$ws.Range("I339").Value = 100108006
$ws.Range("J339").Value = "Plátano"
$ws.Range("K339").Value = "Sin especificar"
$ws.Range("L339").Value = "Primera Maduro"
$ws.Range("M339").Value = 160
$ws.Range("N339").Value = 15000
$ws.Range("O339").Value = 15000
$ws.Range("P339").Value = 15000
$ws.Range("Q339").Value = "`$/caja 20 kilos"
$ws.Range("R339").Value = "Ecuador"
$ws.Range("S339").Value = 750
$ws.Range("T339").Value = 20

# Row 340 - Primera Pintón
$ws.Range("A340").Value = 8
$ws.Range("B340").Value = "Terminal La Palmera de La Serena"
$ws.Range("C340").Value = "Coquimbo"
$ws.Range("D340").Value = 44551
$ws.Range("E340").Value = 4
$ws.Range("F340").Value = "Fruta"
$ws.Range("G340").Value = 100108
$ws.Range("H340").Value = "Tropicales y subtropicales"
$ws.Range("I340").Value = 100108006
$ws.Range("J340").Value = "Plátano"
$ws.Range("K340").Value = "Sin especificar"
$ws.Range("L340").Value = "Primera Pintón"
$ws.Range("M340").Value = 160
$ws.Range("N340").Value = 16000
$ws.Range("O340").Value = 16000
$ws.Range("P340").Value = 16000
$ws.Range("Q340").Value = "`$/caja 20 kilos"
$ws.Range("R340").Value = "Ecuador"
$ws.Range("S340").Value = 800
$ws.Range("T340").Value = 20
